$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws4 = $wb.Worksheets.Item(4)

# --- Remove the old external CRAN hyperlink from I7 (testing sheet) ---
$ws2.Range("I7").Hyperlinks.Delete()
$ws2.Range("I7").ClearContents()

# Give I7, L7 and O7 the "Hyperlink" look (style only, no value/link) like the
# original marker cell used to have.
$ws2.Range("I7").Style = "Hyperlink"
$ws2.Range("L7").Style = "Hyperlink"
$ws2.Range("O7").Style = "Hyperlink"

# --- Shrink the D26:D30 merge down to D26:D27 and clear D28:D30 ---
$ws2.Range("D28:D30").UnMerge()
$ws2.Range("D26:D27").Merge()
$ws2.Range("D28:D30").Clear()

# Seed the shared-string table in the same order the original commit used:
# IrisSample!A1, Irrational_number, Decimal_fractions, IrisSample!A2
$ws2.Range("B33").Value = "IrisSample!A1"
$ws2.Range("B25").Value = "Irrational_number"
$ws2.Range("B28").Value = "Decimal_fractions"
$ws2.Range("B34").Value = "IrisSample!A2"

# --- New hyperlink rows 25 & 26: "Irrational_number" (external, no anchor) ---
$ws2.Hyperlinks.Add($ws2.Range("B25"), "https://en.wikipedia.org/wiki/Irrational_number", "", "", "Irrational_number")
$ws4.Range("G6").Copy()
$ws2.Range("B25").PasteSpecial(-4122)

$ws2.Range("B26").Value = "Irrational_number"
$ws2.Hyperlinks.Add($ws2.Range("B26"), "https://en.wikipedia.org/wiki/Irrational_number", "", "", "Irrational_number")
$ws4.Range("G6").Copy()
$ws2.Range("B26").PasteSpecial(-4122)

# --- New hyperlink rows 28-30: "Decimal_fractions" (external + anchor) ---
$ws2.Hyperlinks.Add($ws2.Range("B28"), "https://en.wikipedia.org/wiki/Irrational_number", "Decimal_fractions", "", "Decimal_fractions")
$ws4.Range("G6").Copy()
$ws2.Range("B28").PasteSpecial(-4122)

$ws2.Range("B29").Value = "Decimal_fractions"
$ws2.Hyperlinks.Add($ws2.Range("B29"), "https://en.wikipedia.org/wiki/Irrational_number", "Decimal_fractions", "", "Decimal_fractions")
$ws4.Range("G6").Copy()
$ws2.Range("B29").PasteSpecial(-4122)

$ws2.Range("B30").Value = "Decimal_fractions"
$ws2.Hyperlinks.Add($ws2.Range("B30"), "https://en.wikipedia.org/wiki/Irrational_number", "Decimal_fractions", "", "Decimal_fractions")
$ws4.Range("G6").Copy()
$ws2.Range("B30").PasteSpecial(-4122)

# --- New internal hyperlink rows 33 & 34, pointing back at IrisSample ---
$ws2.Hyperlinks.Add($ws2.Range("B33"), "", "IrisSample!A1", "", "IrisSample!A1")
$ws2.Range("B33").Style = "Hyperlink"

$ws2.Hyperlinks.Add($ws2.Range("B34"), "", "IrisSample!A2", "", "IrisSample!A2")
$ws2.Range("B34").Style = "Hyperlink"

# --- Make "testing" the active tab (was IrisSample) ---
$ws2.Activate()
